$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the Heading1
#    title paragraph ("Play Arcane Reel Chaos for Free - Review and
#    Features"), with "Meta description" bold and the rest of the
#    sentence in plain text.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$insPoint = $d.Range($metaRange.Start, $metaRange.Start)
$insPoint.InsertAfter("Meta description: Learn how to play and what to expect with Arcane Reel Chaos online slot game for free. Read our review and discover game design, features, and similar games.")

$metaRange2 = $d.Paragraphs(2).Range
$boldRange = $d.Range($metaRange2.Start, $metaRange2.Start + 16)
$boldRange.Bold = 1

# ---------------------------------------------------------------------
# 2) Remove the duplicated bold "Play Arcane Reel Chaos for Free -
#    Review and Features" paragraph near the end of the document (it
#    now duplicates the Heading1 title) and replace the text of the
#    italic paragraph that follows it with the new image-generation
#    prompt, keeping its italic run formatting intact.
# ---------------------------------------------------------------------
$dupFound = $false
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd() -eq "Play Arcane Reel Chaos for Free - Review and Features") {
        $para.Range.Delete()
        $dupFound = $true
        break
    }
}

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$lastTextRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$lastTextRange.Text = "Prompt: Create a feature image for `"Arcane Reel Chaos`" that captures the game's mysterious and supernatural theme, featuring a happy Maya warrior with glasses in cartoon style. The image should have a dark background with the warrior holding a glowing staff or weapon in one hand and a stack of cards in the other, symbolizing the game's use of special cards. The warrior's glasses should reflect the game's logo or a part of the game grid. The overall mood should be adventurous, suggesting a journey into the arcane world of the game."

Write-Output ("edit applied, duplicate title removed=" + $dupFound)
